# Edit: split the "pay period / gross income" bullet into two separate
# list items, and add a new method-focused bullet for gross income.
$d = $word.ActiveDocument

$bullet = [char]0x2022

$found1 = $d.Content.Find.Execute(
    "pay period = per calendar month " + $bullet + " ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pay period = per calendar month ^p",
    2)

# Drop the stale lastRenderedPageBreak marker by re-typing the paragraph's
# text in place (forces the run to be regenerated without the marker).
$found2 = $d.Content.Find.Execute(
    "For example, the payment in March for an employee with an annual salary of `$60,050 and a super rate of 9% is: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "For example, the payment in March for an employee with an annual salary of `$60,050 and a super rate of 9% is: ",
    2)

if (-not $found1) { throw "Could not find the 'pay period / gross income' bullet to split." }
if (-not $found2) { throw "Could not find the 'For example, the payment in March...' sentence to refresh." }

Write-Output "split: $found1; repaginate: $found2"
